$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Header text updates (shared strings used by both sheets) ---
$ws1.Range("A3").Value = "Precincts 307 of 307 Reporting (Precincts Partially Reported: 0/307)"
$ws1.Range("A4").Value = "Downloaded at 7/19/2022 9:10:38 AM"
$ws2.Range("A3").Value = "Precincts 307 of 307 Reporting (Precincts Partially Reported: 0/307)"
$ws2.Range("A4").Value = "Downloaded at 7/19/2022 9:10:38 AM"

# --- Column A width update (both sheets) ---
$ws1.Columns.Item(1).ColumnWidth = 58.5
$ws2.Columns.Item(1).ColumnWidth = 58.5

# --- Sheet 1 (SUPREME COURT JUSTICE #1) vote count updates ---
$ws1.Range("C8").Value = 1863
$ws1.Range("D8").Value = 432
$ws1.Range("C11").Value = 1333
$ws1.Range("C12").Value = 2137
$ws1.Range("C14").Value = 12574
$ws1.Range("D14").Value = 3400
$ws1.Range("C16").Value = 2186
$ws1.Range("D16").Value = 641
$ws1.Range("C18").Value = 1708
$ws1.Range("C19").Value = 1738
$ws1.Range("D19").Value = 692
$ws1.Range("C21").Value = 2673
$ws1.Range("D21").Value = 642
$ws1.Range("C22").Value = 18017
$ws1.Range("D22").Value = 5377
$ws1.Range("C23").Value = 17519
$ws1.Range("D23").Value = 5708
$ws1.Range("C25").Value = 1303
$ws1.Range("D25").Value = 467
$ws1.Range("C27").Value = 754
$ws1.Range("C28").Value = 2070
$ws1.Range("D28").Value = 555
$ws1.Range("C29").Value = 3078
$ws1.Range("D29").Value = 650
$ws1.Range("C31").Value = 5153
$ws1.Range("D31").Value = 1702
$ws1.Range("C32").Value = 15187
$ws1.Range("D32").Value = 3385
$ws1.Range("C34").Value = 3375
$ws1.Range("D34").Value = 1227
$ws1.Range("C35").Value = 2019
$ws1.Range("D35").Value = 603
$ws1.Range("D38").Value = 326
$ws1.Range("C39").Value = 20115
$ws1.Range("D39").Value = 7791
$ws1.Range("C41").Value = 3107
$ws1.Range("D41").Value = 1129
$ws1.Range("C44").Value = 1151
$ws1.Range("D44").Value = 293
$ws1.Range("C48").Value = 8911
$ws1.Range("D48").Value = 2855
$ws1.Range("C49").Value = 1557
$ws1.Range("D49").Value = 458
$ws1.Range("C51").Value = 1163
$ws1.Range("D51").Value = 371
$ws1.Range("C52").Value = 3012
$ws1.Range("D52").Value = 933
$ws1.Range("C53").Value = 696
$ws1.Range("C54").Value = 5954
$ws1.Range("D54").Value = 1781
$ws1.Range("D55").Value = 504
$ws1.Range("C56").Value = 997
$ws1.Range("D56").Value = 207
$ws1.Range("C58").Value = 820
$ws1.Range("D58").Value = 251
$ws1.Range("D61").Value = 104
$ws1.Range("C63").Value = 28985
$ws1.Range("D63").Value = 10844
$ws1.Range("C64").Value = 189101
$ws1.Range("D64").Value = 59168

# --- Sheet 2 (SUPREME COURT JUSTICE #2) vote count updates ---
$ws2.Range("C8").Value = 870
$ws2.Range("D8").Value = 318
$ws2.Range("E8").Value = 1295
$ws2.Range("C11").Value = 519
$ws2.Range("E11").Value = 747
$ws2.Range("C14").Value = 7555
$ws2.Range("D14").Value = 3271
$ws2.Range("E14").Value = 5593
$ws2.Range("C16").Value = 1121
$ws2.Range("D16").Value = 470
$ws2.Range("E16").Value = 1346
$ws2.Range("C18").Value = 807
$ws2.Range("D18").Value = 321
$ws2.Range("E18").Value = 880
$ws2.Range("C19").Value = 1349
$ws2.Range("D19").Value = 631
$ws2.Range("E19").Value = 543
$ws2.Range("C21").Value = 1007
$ws2.Range("D21").Value = 575
$ws2.Range("E21").Value = 1867
$ws2.Range("C22").Value = 9874
$ws2.Range("D22").Value = 3594
$ws2.Range("E22").Value = 11476
$ws2.Range("C23").Value = 14200
$ws2.Range("D23").Value = 2228
$ws2.Range("E23").Value = 8403
$ws2.Range("C25").Value = 1001
$ws2.Range("D25").Value = 304
$ws2.Range("E25").Value = 559
$ws2.Range("C28").Value = 1152
$ws2.Range("D28").Value = 527
$ws2.Range("E28").Value = 984
$ws2.Range("C29").Value = 1583
$ws2.Range("D29").Value = 889
$ws2.Range("E29").Value = 1355
$ws2.Range("C31").Value = 3209
$ws2.Range("D31").Value = 1108
$ws2.Range("E31").Value = 2995
$ws2.Range("C32").Value = 8818
$ws2.Range("D32").Value = 6424
$ws2.Range("E32").Value = 4181
$ws2.Range("C34").Value = 1743
$ws2.Range("D34").Value = 960
$ws2.Range("E34").Value = 2108
$ws2.Range("C35").Value = 1298
$ws2.Range("D35").Value = 399
$ws2.Range("E35").Value = 1074
$ws2.Range("D38").Value = 217
$ws2.Range("C39").Value = 19241
$ws2.Range("D39").Value = 3527
$ws2.Range("E39").Value = 7094
$ws2.Range("C41").Value = 2315
$ws2.Range("D41").Value = 415
$ws2.Range("E41").Value = 1785
$ws2.Range("C44").Value = 626
$ws2.Range("D44").Value = 253
$ws2.Range("E44").Value = 610
$ws2.Range("D45").Value = 83
$ws2.Range("C48").Value = 4947
$ws2.Range("D48").Value = 1751
$ws2.Range("E48").Value = 5793
$ws2.Range("C49").Value = 713
$ws2.Range("D49").Value = 480
$ws2.Range("E49").Value = 834
$ws2.Range("C51").Value = 644
$ws2.Range("D51").Value = 215
$ws2.Range("E51").Value = 759
$ws2.Range("C52").Value = 1287
$ws2.Range("D52").Value = 594
$ws2.Range("E52").Value = 2266
$ws2.Range("C53").Value = 394
$ws2.Range("D53").Value = 179
$ws2.Range("C54").Value = 4611
$ws2.Range("D54").Value = 1808
$ws2.Range("E54").Value = 1690
$ws2.Range("C55").Value = 1059
$ws2.Range("E55").Value = 1193
$ws2.Range("C56").Value = 511
$ws2.Range("D56").Value = 158
$ws2.Range("E56").Value = 651
$ws2.Range("C58").Value = 361
$ws2.Range("D58").Value = 222
$ws2.Range("E58").Value = 528
$ws2.Range("C63").Value = 24207
$ws2.Range("D63").Value = 4413
$ws2.Range("E63").Value = 14578
$ws2.Range("C64").Value = 126423
$ws2.Range("D64").Value = 40872
$ws2.Range("E64").Value = 95607
